$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'mens compression basketball pants',
    'mens basketball tights with knee pads',
    'capri pants for boys',
    'knee pads with leggings',
    'compression pants for men',
    'compression mens leggings',
    'basketball protective gear',
    'basketball spandex pants',
    'compression tights girls',
    'wrestling pants',
    'mens compression capris',
    'baseball sliding tights',
    'boys basketball compression pants with knee pads',
    'compression men capri',
    'volleyball knee pads men',
    'mens basketball',
    'athletic leggings boys',
    'knee pad weightlifting',
    'basketball guide',
    'basketball knee pads boys',
    'youth padded compression pants basketball',
    'leggings for men',
    'cycling capris',
    'padded knee pads for basketball',
    'pant with knee pad',
    'tight pants for men',
    'basketball kneepads',
    'thigh compression leggings',
    'compression pants men pack',
    'hiking knee pads',
    'sliding pants',
    'boys compression pants',
    'cycling knee pads for men',
    'youth boys knee pads for basketball',
    'mens spandex leggings',
    'compression knee padded',
    'youth basketball kneepads',
    'baseball equipment for boys',
    'youth football tights boys',
    'men compression workout pants',
    'padded knee basketball',
    'mens compression pants running',
    'basketball youth leggings',
    'compression basketball leggings youth',
    'men compression running pants',
    'wrestling knee pads youth',
    'mens leggings sports',
    'skin leggings men',
    'youth football pants',
    'mens leggings for sports',
    'black boys compression pants',
    'basketball team clothes',
    'men compression tights',
    'youth spandex leggings boys',
    'baseball gear for men',
    'sports knee pads',
    'leggings men',
    'pain in hip down leg',
    'compression men tights',
    'softball pants for men',
    'youth basketball knee pads for boys',
    'mens running tights',
    'youth baseball gear',
    'football pants adult xl',
    'padded knee compression',
    'compression pants men 3/4',
    'sliding pants softball youth',
    'hockey hip pads',
    'compression knees',
    'compression men',
    'compression mens tights',
    'compression tights men',
    'legging for men',
    'football knee pads adult',
    'football pants with pads mens',
    'mens basketball tights and leggings',
    'knee protector',
    'basketball apparel mens',
    'adult knee pads',
    'black softball pants girls',
    'youth football pads for pants',
    'youth black compression pants',
    '3/4 compression pants men',
    'girls hiking pants',
    'youth tights boys basketball',
    'down pants men',
    'basketball volleyball knee pads',
    'legging for men sport',
    'polyester hex mesh',
    'basketball knee pad tights',
    'capris tights',
    'thigh compression pants',
    'bjj pants men',
    'wrestling pads',
    'basketball athletic tights',
    'wrestling knee pads youth 2 pack',
    'compression for knees',
    'youth football pants with pads black',
    'girls knee pads',
    'compression pants for work'
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
